# gsc-export/HTTPS.xlsx -- append the next day's row (2025-12-24) to the
# "Chart" sheet's data table, same as the daily GSC export refresh.

$wb = $excel.ActiveWorkbook
$chart = $wb.Worksheets.Item("Chart")

$newRow = $chart.UsedRange.Rows.Count + 1

# Column A holds the date as plain text (e.g. "2025-10-07"). A leading
# apostrophe forces Excel to store it as text instead of auto-converting
# the recognizable yyyy-mm-dd pattern into a date serial number.
$chart.Cells.Item($newRow, 1).Value = "'2025-12-24"
$chart.Cells.Item($newRow, 1).ClearFormats()

$chart.Cells.Item($newRow, 2).Value = 0
$chart.Cells.Item($newRow, 3).Value = 30
